# "Generate Report for Handback"
# Updates the "Latest HO Xliff Generate Date" on the Overview sheet and the
# corresponding "Correspond Handoff/Handback DateTime" columns on the
# per-locale (zh-cn / de-de) sheets for the ecf43ef3-... row, reflecting a
# fresh handback run.

$wb = $excel.ActiveWorkbook

# Overview sheet: row for ecf43ef3-1d62-4099-80e8-e264b6095bb6.md, de-de column (G)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-22 04:57:18"

# zh-cn sheet: row 3 is the ecf43ef3-... file
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-08-22 04:57:13"
$wsZhCn.Range("K3").Value = "2016-08-22 04:57:31"

# de-de sheet: row 3 is the ecf43ef3-... file
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-08-22 04:57:18"
$wsDeDe.Range("K3").Value = "2016-08-22 04:57:37"
